# تعديل يدوي في شيت Card19 by admin at 2025-12-17 13:44:30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

# D8 holds a numeric-looking value that must stay stored as text (like the
# rest of the sheet's "nan"/number-as-text columns), so force text format
# before entry and drop the number-format afterwards so the cell keeps the
# sheet's default (unstyled) look.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "993"
$ws.Range("D8").ClearFormats()

$ws.Range("F8").Value = "✅"
$ws.Range("K8").Value = "✅"
$ws.Range("L8").Value = "17/8/2025"
$ws.Range("O8").Value = "م.محمد عبدالله ،ف.محمود ايهاب ،حسام،عمر"
